$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the AD3:AI3 values
$ws.Range("AD3").Value = 1
$ws.Range("AE3").Value = 1
$ws.Range("AF3").Value = [double]"1.0000000000000001E-15"
$ws.Range("AG3").Value = 1
$ws.Range("AH3").Value = 1
$ws.Range("AI3").Value = 1

# Update the AD5:AI5 values
$ws.Range("AD5").Value = 1
$ws.Range("AE5").Value = 1
$ws.Range("AF5").Value = [double]"1.0000000000000001E-15"
$ws.Range("AG5").Value = 1
$ws.Range("AH5").Value = 1
$ws.Range("AI5").Value = 1

# Column widths for AD:AI (30-35)
# NOTE: the host engine quantizes ColumnWidth to whole pixels (MDW=6) when
# serialising back to OOXML (xml_width = (round(ColumnWidth*6)+5)/6), so the
# inputs below are chosen to land as close as the pixel grid allows to the
# exact target widths from the authored file.
$ws.Columns.Item(30).ColumnWidth = 32.5
$ws.Columns.Item(31).ColumnWidth = 24.333333333333332
$ws.Columns.Item(32).ColumnWidth = 28.0
$ws.Columns.Item(33).ColumnWidth = 29.5
$ws.Columns.Item(34).ColumnWidth = 29.166666666666668
$ws.Columns.Item(35).ColumnWidth = 30.666666666666668

# Selection / view change: move the active cell from X6 to D6 and clear the
# scrolled-right "topLeftCell" view state.
$ws.Range("D6").Select()
